$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("composite")
$ws.Columns.Item(7).AutoFit() | Out-Null
$ws.Columns.Item(8).AutoFit() | Out-Null
Write-Host "G width:" $ws.Columns.Item(7).ColumnWidth
Write-Host "H width:" $ws.Columns.Item(8).ColumnWidth
